$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.018.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.14%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.352.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.16%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.12%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'545.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.93%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'137.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.21%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.15%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -8.50%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.353.12"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.07%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.18%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +1.47%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'5.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.55%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.59%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'24.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.23%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.777.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.00%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'60.863.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.02%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -2.07%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.349.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.16%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'10.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.36%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'320.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.24%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +1.05%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.15%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.05%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'1.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -4.39%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.60%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'8.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +9.20%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.28%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.467.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.03%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.45%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'501.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.56%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -2.73%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.0₃0871"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -6.50%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +2.31%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -1.62%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -3.52%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.05%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.21%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +1.32%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +2.97%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'5.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.06%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +6.72%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'142.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +3.74%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'40.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.64%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'142.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +2.33%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +1.23%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -5.87%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +0.49%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'19.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -5.24%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.570"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.72%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -1.96%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0221"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.79%  "
$ws.Range("E51").Style = "Normal"
